$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.560471058
$ws.Range("C2").Value = -224.73773952
$ws.Range("D2").Value = -225.29821058

$ws.Range("B3").Value = -0.5691579122
$ws.Range("C3").Value = -224.66015071
$ws.Range("D3").Value = -225.22930862

$ws.Range("B4").Value = -0.5730019284
$ws.Range("C4").Value = -224.64730749
$ws.Range("D4").Value = -225.22030942
